$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 3
$ws.Range("E3").Value = 4
$ws.Range("G3").Value = -3
$ws.Range("H3").Value = 13

# Select cell D3 (reflected in sheetView selection)
$ws.Range("D3").Select()
